# Inserts two new data rows (new rows 175 and 176) into the "Frambuesa"
# sheet, pushing the former rows 175-219 down to 177-221.
#
# The new rows carry the same constant columns (A,B,C,E,F,G,H,I,J,K,Q) that
# every data row in this sheet shares, with the following specific values:
#   Row 175: D=44900 L=Especial M=250 N=8000 O=8000 P=8000 R="Región de O'Higgins" S=4000 T=2
#   Row 176: D=44900 L=Especial M=350 N=8000 O=8000 P=8000 R="Región del Maule"   S=4000 T=2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 175:221 down by inserting two blank rows at row 175.
$ws.Rows.Item(175).Resize(2).Insert()

# Constant values shared by every data row in this sheet.
$constA = 6
$constB = "Mercado Mayorista Lo Valledor de Santiago"
$constC = "Metropolitana"
$constE = 13
$constF = "Fruta"
$constG = 100101
$constH = "Berries"
$constI = 100101004
$constJ = "Frambuesa"
$constK = "Sin especificar"
$constQ = "`$/bandeja 2 kilos"

function Set-DataRow {
    param($Row, $D, $L, $M, $N, $O, $P, $R, $S, $T)

    $ws.Cells.Item($Row, 1).Value = $constA
    $ws.Cells.Item($Row, 2).Value = $constB
    $ws.Cells.Item($Row, 3).Value = $constC
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $constE
    $ws.Cells.Item($Row, 6).Value = $constF
    $ws.Cells.Item($Row, 7).Value = $constG
    $ws.Cells.Item($Row, 8).Value = $constH
    $ws.Cells.Item($Row, 9).Value = $constI
    $ws.Cells.Item($Row, 10).Value = $constJ
    $ws.Cells.Item($Row, 11).Value = $constK
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $constQ
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 19).Value = $S
    $ws.Cells.Item($Row, 20).Value = $T
}

Set-DataRow 175 44900 "Especial" 250 8000 8000 8000 "Región de O'Higgins" 4000 2
Set-DataRow 176 44900 "Especial" 350 8000 8000 8000 "Región del Maule"    4000 2

# Match the date number format already used by the other rows in column D.
$ws.Range("D175:D176").NumberFormat = $ws.Range("D174").NumberFormat
